# Update workbook to reflect new data pulled through 2021-10-13
# (commit message: "Add data for 2021-10-21")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-13"

# --- Row 9 (July) : 2021 arrest/no_arrest/arrest_rate updated ---
$ws.Range("T9").Value = 12
$ws.Range("U9").Value = 139
$ws.Range("V9").Value = 0.0795

# --- Row 12 (October, through date) ---
$ws.Range("A12").Value = "October (through 10-13)"

$ws.Range("E12").Value = 1
$ws.Range("G12").Value = 0.0526

$ws.Range("H12").Value = 5
$ws.Range("I12").Value = 16
$ws.Range("J12").Value = 0.2381

$ws.Range("L12").Value = 30
$ws.Range("M12").Value = 0.0625

$ws.Range("O12").Value = 16
$ws.Range("P12").Value = 0.0588

$ws.Range("R12").Value = 63
$ws.Range("U12").Value = 80

# --- Row 13 (Total) ---
$ws.Range("E13").Value = 47
$ws.Range("G13").Value = 0.1049

$ws.Range("H13").Value = 55
$ws.Range("I13").Value = 593
$ws.Range("J13").Value = 0.0849

$ws.Range("L13").Value = 517
$ws.Range("M13").Value = 0.1086

$ws.Range("O13").Value = 395
$ws.Range("P13").Value = 0.1002

$ws.Range("R13").Value = 911
$ws.Range("S13").Value = 0.055

$ws.Range("T13").Value = 81
$ws.Range("V13").Value = 0.0609
